$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Delete the whole "Corrected networking connectivity issues ..." bullet
#    paragraph (list item under numId 37), including its own paragraph mark,
#    so the following bullet ("Router Programming and automating ...")
#    becomes adjacent to the "... and then parse to analyze different
#    results." paragraph.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    if ($txt -like "*Corrected networking connectivity issues in wireless, routing, and switching using a layered model approach*") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark (Word's "last edit location" bookmark) from
#    the end of the "Design, code, and test networking ..." paragraph to the
#    end of the "... and then parse to analyze different results." paragraph
#    - i.e. remove it from its old spot and re-insert it, collapsed, right
#    after the last run of the new target paragraph (before its paragraph
#    mark), exactly mirroring where Word leaves it after the edit above.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text
    if ($txt -like "*and then parse to analyze different results*") {
        # Position right before this paragraph's own paragraph mark.
        $endPos = $p.Range.End - 1

        # A zero-length Range sitting exactly on a paragraph-mark boundary
        # cannot be handed to Bookmarks.Add directly and reliably resolve to
        # that boundary, so: insert a throwaway marker character there,
        # wrap the bookmark around it, then delete the marker. Deleting the
        # bracketed character collapses the bookmark back down to a
        # zero-width bookmark sitting exactly where the marker was, which is
        # precisely the position we want.
        $d.Range($endPos, $endPos).InsertAfter("X")
        $d.Bookmarks.Add("_GoBack", $d.Range($endPos, $endPos + 1))
        $d.Range($endPos, $endPos + 1).Delete()
        break
    }
}
